# Commit: "Added more options for creation of design matrix, and senstype ref."
#
# 1. Rename the parameter headers in row 1 to their upper-case form
#    (param13/14/15/16 -> PARAM13/14/15/16).
# 2. Narrow the active selection on Sheet1 from D2:D12 down to just D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "PARAM13"
$ws.Range("B1").Value = "PARAM14"
$ws.Range("C1").Value = "PARAM15"
$ws.Range("D1").Value = "PARAM16"

[void]$ws.Range("D2").Select()
